$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sourcing")

# Row 12 - "Connecteur d'alimentation": diameter to be determined
$ws.Range("G12").Value = "Diamètre à définir"

# Row 9 - "Porte fusible": add alternative reference (G9) and price (H9)
$ws.Range("G9").Value = "BK1/HTC-100M"
$ws.Range("H9").Value = 0.722

# "Nb PCB" label (A15) - nudge the alignment off and back onto
# center/center so the cell's cached style entry gets rewritten without
# the stale (no-op) fill/border "apply" flags it carried before.
$lbl = $ws.Range("A15")
$lbl.HorizontalAlignment = -4131
$lbl.VerticalAlignment = -4160
$lbl.HorizontalAlignment = -4108
$lbl.VerticalAlignment = -4108

# Update active selection to reflect where the user ended up
$ws.Range("G16").Select()
